$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tables")
$ws.Activate()

$ws.Range("B3").Value = "A"
$ws.Range("B4").Value = "Z"
$ws.Range("G4").Formula = '=MATCH("ZB",FirstTable[#This Row],FALSE)'
$ws.Range("G17").Formula = "=IF(C4>D4,TRUE)"

$ws.Range("G5").Select() | Out-Null
